$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the Kathmandu row (currently row 38: "Kathmandu, Lalitpur Nepal" / "Kathmandu")
$valA = $ws.Cells.Item(38, 1).Value2
$valB = $ws.Cells.Item(38, 2).Value2

# Shift rows 25-37 down into rows 26-38 (iterate bottom-up to avoid overwriting).
for ($r = 37; $r -ge 25; $r--) {
    $ws.Cells.Item($r + 1, 1).Value2 = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r + 1, 2).Value2 = $ws.Cells.Item($r, 2).Value2
}

# Place the Kathmandu row into its new alphabetical position at row 25.
$ws.Cells.Item(25, 1).Value2 = $valA
$ws.Cells.Item(25, 2).Value2 = $valB

# Restore/update the view state: clear the scrolled top-left cell and set the
# active selection to A6 as recorded after the edit.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("A6").Select()
